# Added bar codes for Irminger D3
$wb = $excel.ActiveWorkbook

$wsMoorings = $wb.Worksheets.Item("Moorings")
$wsAsset = $wb.Worksheets.Item("Asset_Cal_Info")

# Fill in the new barcodes in the order they were keyed in, so the
# shared-string table grows in the same sequence as the authored workbook.

# ENG row (11) first: sensor barcode, then mooring barcode
$wsAsset.Range("E11").Value = "OL000119"
$wsAsset.Range("B11").Value = "A00966"

# FLORD rows (2-5): mooring barcode already known, sensor barcode is new
$wsAsset.Range("B2").Value = "A00966"
$wsAsset.Range("E2").Value = "N00530"
$wsAsset.Range("B3").Value = "A00966"
$wsAsset.Range("E3").Value = "N00530"
$wsAsset.Range("B4").Value = "A00966"
$wsAsset.Range("E4").Value = "N00530"
$wsAsset.Range("B5").Value = "A00966"
$wsAsset.Range("E5").Value = "N00530"

# DOSTA row (7)
$wsAsset.Range("B7").Value = "A00966"
$wsAsset.Range("E7").Value = "N00529"

# CTDGV row (9)
$wsAsset.Range("B9").Value = "A00966"
$wsAsset.Range("E9").Value = "N00528"

# Moorings!A2 - Mooring OOIBARCODE (string already exists by now)
$wsMoorings.Range("A2").Value = "A00966"

# Update selections / active sheet to match the authored state
$wsMoorings.Range("D20").Select() | Out-Null
$wsAsset.Select() | Out-Null
$wsAsset.Range("C22").Select() | Out-Null
